$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 895.0909
$ws.Range("J9").Value = 2283
$ws.Range("N9").Value = -2621
$ws.Range("L9").Value = 2283
$ws.Range("K9").Value = 374.625
$ws.Range("M9").Value = -205.625
$ws.Range("I9").Value = 374.625
$ws.Range("H40").Value = 7021.1113
$ws.Range("J40").Value = 7862
$ws.Range("N40").Value = -8212
$ws.Range("L40").Value = 7862
$ws.Range("N43").Value = -3159.9092
$ws.Range("H43").Value = 5032.4165
$ws.Range("M43").Value = -5848.04
$ws.Range("K43").Value = 5917.04
$ws.Range("I43").Value = 5917.04
$ws.Range("L43").Value = 3021.9092
$ws.Range("J43").Value = 3021.9092
$ws.Range("H106").Value = 48891696
$ws.Range("K106").Value = 62859796
$ws.Range("I106").Value = 62859796
$ws.Range("M106").Value = -62859165
$ws.Range("K107").Value = 2003
$ws.Range("I107").Value = 2003
$ws.Range("H107").Value = 2375.875
$ws.Range("M107").Value = -83

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I110").Value = 4888.8887
$ws.Range("H110").Value = 4888.8887
$ws.Range("M110").Value = -2843.8887
$ws.Range("K110").Value = 4888.8887
$ws.Range("J133").Value = 38807.94
$ws.Range("N133").Value = -43867.94
$ws.Range("H133").Value = 38807.94
$ws.Range("L133").Value = 38807.94

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N86").Value = -3942.4445
$ws.Range("M86").Value = -1130.8
$ws.Range("I86").Value = 2253.8
$ws.Range("L86").Value = 1696.4445
$ws.Range("J86").Value = 1696.4445
$ws.Range("K86").Value = 2253.8
$ws.Range("H86").Value = 1989.7894
$ws.Range("I89").Value = 2253.8
$ws.Range("M89").Value = -5653
$ws.Range("L89").Value = 8482.2225
$ws.Range("K89").Value = 11269
$ws.Range("H89").Value = 1989.7894
$ws.Range("N89").Value = -19714.2225
$ws.Range("J89").Value = 1696.4445
$ws.Range("L100").Value = 17492.5
$ws.Range("H100").Value = 17492.5
$ws.Range("J100").Value = 17492.5
$ws.Range("N100").Value = -19656.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("M16").Value = -374
$ws.Range("K16").Value = 661
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("I16").Value = 661
$ws.Range("H16").Value = 661
$ws.Range("J113").Value = 0
$ws.Range("I113").Value = 661
$ws.Range("H113").Value = 661
$ws.Range("M113").Value = 1509
$ws.Range("K113").Value = 661
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("I134").Value = 166670990
$ws.Range("J134").Value = 7999.6665
$ws.Range("L134").Value = 23998.9995
$ws.Range("N134").Value = -29068.9995
$ws.Range("H134").Value = 83339496
$ws.Range("M134").Value = -500010435
$ws.Range("K134").Value = 500012970
$ws.Range("L141").Value = 227857
$ws.Range("H141").Value = 227857
$ws.Range("N141").Value = -238217
$ws.Range("J141").Value = 227857

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N22").ClearContents()
$ws.Range("J22").Value = 0
$ws.Range("H22").Value = 1950
$ws.Range("L22").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("J27").Value = 0
$ws.Range("H27").Value = 1950
$ws.Range("L27").Value = 0
$ws.Range("M68").Value = -2063.75
$ws.Range("I68").Value = 958.25
$ws.Range("L68").Value = 9009
$ws.Range("J68").Value = 3003
$ws.Range("N68").Value = -10631
$ws.Range("H68").Value = 1367.2
$ws.Range("K68").Value = 2874.75
$ws.Range("N71").Value = -35139
$ws.Range("M71").Value = -4568.25
$ws.Range("J71").Value = 3003
$ws.Range("H71").Value = 1367.2
$ws.Range("L71").Value = 27027
$ws.Range("I71").Value = 958.25
$ws.Range("K71").Value = 8624.25
$ws.Range("H99").Value = 1999.3334
$ws.Range("J99").Value = 5000
$ws.Range("N99").Value = -19492
$ws.Range("L99").Value = 15000
$ws.Range("K107").Value = 735
$ws.Range("I107").Value = 245
$ws.Range("J107").Value = 620.9091
$ws.Range("N107").Value = -5702.7273
$ws.Range("H107").Value = 503.4375
$ws.Range("L107").Value = 1862.7273
$ws.Range("M107").Value = 1185
$ws.Range("J113").Value = 4348591.5
$ws.Range("I113").Value = 475
$ws.Range("H113").Value = 3623905.5
$ws.Range("M113").Value = 745
$ws.Range("K113").Value = 1425
$ws.Range("L113").Value = 13045774.5
$ws.Range("N113").Value = -13050114.5
$ws.Range("H118").Value = 2202.5
$ws.Range("M118").Value = -5364.5
$ws.Range("I118").Value = 2202.5
$ws.Range("K118").Value = 6607.5
$ws.Range("L121").Value = 3799637.4
$ws.Range("N121").Value = -3802257.4
$ws.Range("J121").Value = 1266545.8
$ws.Range("H121").Value = 905803.4

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I102").Value = 4198.5713
$ws.Range("H102").Value = 5297.643
$ws.Range("K102").Value = 4198.5713
$ws.Range("M102").Value = -2576.5713
$ws.Range("J113").Value = 1648.1818
$ws.Range("I113").Value = 1485.88
$ws.Range("H113").Value = 1535.4722
$ws.Range("M113").Value = 684.1199999999999
$ws.Range("K113").Value = 1485.88
$ws.Range("L113").Value = 1648.1818
$ws.Range("N113").Value = -5988.1818
$ws.Range("H122").Value = 2486.8635
$ws.Range("J122").Value = 6659.6
$ws.Range("K122").Value = 3778.7649
$ws.Range("M122").Value = -1328.7649
$ws.Range("I122").Value = 1259.5883
$ws.Range("L122").Value = 19978.8
$ws.Range("N122").Value = -24878.8
$ws.Range("J140").Value = 99200
$ws.Range("L140").Value = 99200
$ws.Range("N140").Value = -109560
$ws.Range("H140").Value = 99200

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J19").Value = 49200
$ws.Range("N19").Value = -49540
$ws.Range("H19").Value = 49200
$ws.Range("L19").Value = 49200
$ws.Range("M46").Value = -659.9091
$ws.Range("K46").Value = 847.9091
$ws.Range("H46").Value = 2697.85
$ws.Range("I46").Value = 847.9091
$ws.Range("M68").Value = -5514.1665
$ws.Range("I68").Value = 6263.1665
$ws.Range("L68").Value = 4699.6665
$ws.Range("J68").Value = 4699.6665
$ws.Range("N68").Value = -6197.6665
$ws.Range("H68").Value = 5950.467
$ws.Range("K68").Value = 6263.1665
$ws.Range("N71").Value = -30986.3325
$ws.Range("M71").Value = -27571.8325
$ws.Range("J71").Value = 4699.6665
$ws.Range("H71").Value = 5950.467
$ws.Range("L71").Value = 23498.3325
$ws.Range("I71").Value = 6263.1665
$ws.Range("K71").Value = 31315.8325
$ws.Range("H122").Value = 9462.125
$ws.Range("K122").Value = 9999
$ws.Range("M122").Value = -7549
$ws.Range("I122").Value = 3333
$ws.Range("J133").Value = 99531.5
$ws.Range("N133").Value = -104591.5
$ws.Range("H133").Value = 99531.5
$ws.Range("L133").Value = 99531.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K62").Value = 6999.5
$ws.Range("H62").Value = 9349.799999999999
$ws.Range("I62").Value = 6999.5
$ws.Range("M62").Value = -6375.5
$ws.Range("K65").Value = 34997.5
$ws.Range("H65").Value = 9349.799999999999
$ws.Range("M65").Value = -31877.5
$ws.Range("I65").Value = 6999.5
$ws.Range("I113").Value = 488.5
$ws.Range("H113").Value = 514.5454999999999
$ws.Range("M113").Value = 704.5
$ws.Range("K113").Value = 1465.5
